$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AZ (52), shifting the existing "Mean" column (AZ -> BA)
$ws.Columns.Item(52).Insert()

# New column header: "Run 50"
$ws.Range("AZ1").Value = "Run 50"

# Fill the new "Run 50" column with its value for every data row (rows 2-14)
$azValue = 80.77347764
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = $azValue
}

# Update the "Mean" column (now BA) with the recalculated mean that includes Run 50
$meanValue = 116.24383204
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 53).Value = $meanValue
}
